$d = $word.ActiveDocument
$q = [char]34

# ---------------------------------------------------------------------
# 1) The three pictures were re-inserted/refreshed, which stamps each
#    containing run with <w:rPr><w:noProof/></w:rPr> so the proofer
#    skips the (non-text) drawing run. Mirror that on every InlineShape.
# ---------------------------------------------------------------------
foreach ($ishp in $d.InlineShapes) {
    $ishp.Range.NoProofing = $true
}

# ---------------------------------------------------------------------
# 2) In the VBA listing, the line that sets the NumberFormat of column L
#    (". Range("L" & SumRecord - 1).NumberFormat") had its "gramStart" /
#    "gramEnd" grammar-check markers removed (the "-1.Range" run and the
#    following ("L" & run got merged into a single run), and the format
#    string changed from "#,##0.00" to "#,##0".
#    Locate the paragraph by its (unique) text rather than a fixed
#    offset, so the edit is resilient to any earlier content shifting
#    the story's character positions.
# ---------------------------------------------------------------------
$target = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like ("*SumRecord - 1).NumberFormat = " + $q + "#,##0.00" + $q + "*")) {
        $target = $p.Range
    }
}

if ($target -ne $null) {
    # Exclude the trailing paragraph-mark character from the range we
    # clear/replace so only the paragraph's content changes.
    $p0 = $target.Start
    $p1 = $target.End - 1
    $rng = $d.Range($p0, $p1)

    $rng.Text = ""

    $newPara = "<w:p w14:paraId=" + $q + "1809F2CD" + $q + " w14:textId=" + $q + "77777777" + $q + " w:rsidR=" + $q + "00734557" + $q + " w:rsidRDefault=" + $q + "00734557" + $q + " w:rsidP=" + $q + "00734557" + $q + ">" + `
        "<w:r><w:t xml:space=" + $q + "preserve" + $q + ">       Worksheets(</w:t></w:r>" + `
        "<w:proofErr w:type=" + $q + "spellStart" + $q + "/><w:r><w:t>ws</w:t></w:r><w:proofErr w:type=" + $q + "spellEnd" + $q + "/>" + `
        "<w:r><w:t xml:space=" + $q + "preserve" + $q + ">).Range(" + $q + "L" + $q + " &amp; </w:t></w:r>" + `
        "<w:proofErr w:type=" + $q + "spellStart" + $q + "/><w:r><w:t>SumRecord</w:t></w:r><w:proofErr w:type=" + $q + "spellEnd" + $q + "/>" + `
        "<w:r><w:t xml:space=" + $q + "preserve" + $q + "> - 1).</w:t></w:r>" + `
        "<w:proofErr w:type=" + $q + "spellStart" + $q + "/><w:r><w:t>NumberFormat</w:t></w:r><w:proofErr w:type=" + $q + "spellEnd" + $q + "/>" + `
        "<w:r><w:t xml:space=" + $q + "preserve" + $q + "> = " + $q + "#,##0" + $q + "</w:t></w:r>" + `
        "</w:p>"

    $xmlFrag = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
        '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + `
        $newPara + `
        '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

    $rng.InsertXML($xmlFrag)
}
